$wb = $excel.ActiveWorkbook

# --- Step 1: create the new "2022-Q4" sheet by copying "2022-Q3" (keeps header/style layout) ---
$srcQ3 = $wb.Worksheets.Item("2022-Q3")
$firstSheet = $wb.Worksheets.Item(1)
$srcQ3.Copy($null, $firstSheet)
$newQ4 = $wb.Worksheets.Item(2)
$newQ4.Name = "2022-Q4"

# extend from 33 data rows (Q3) to 39 data rows (Q4): copy style of column A down to the new rows
$newQ4.Range("A34").Copy($newQ4.Range("A35:A40"))

# --- Step 2: write the 2022-Q4 header row (row 1) ---
$newQ4.Range("B1").Value = "基金代码"
$newQ4.Range("C1").Value = "基金名称"
$newQ4.Range("D1").Value = "基金规模"
$newQ4.Range("E1").Value = "股票总仓位"
$newQ4.Range("F1").Value = "仓位占比"
$newQ4.Range("G1").Value = "持有市值(亿元)"
$newQ4.Range("H1").Value = "仓位排名"

# --- Step 3: write the 2022-Q4 data rows (rows 2..40) ---
$newQ4.Range("A2").Value = 0
$newQ4.Range("B2").Value = "007130"
$newQ4.Range("C2").Value = "中庚小盘价值股票"
$newQ4.Range("D2").Value = "76.15"
$newQ4.Range("E2").Value = "93.50"
$newQ4.Range("F2").Value = "5.03"
$newQ4.Range("G2").Value = "3.8303"
$newQ4.Range("H2").Value = 4
$newQ4.Range("A3").Value = 1
$newQ4.Range("B3").Value = "008638"
$newQ4.Range("C3").Value = "广发科技创新混合A"
$newQ4.Range("D3").Value = "33.72"
$newQ4.Range("E3").Value = "91.92"
$newQ4.Range("F3").Value = "6.45"
$newQ4.Range("G3").Value = "2.1749"
$newQ4.Range("H3").Value = 3
$newQ4.Range("A4").Value = 2
$newQ4.Range("B4").Value = "007497"
$newQ4.Range("C4").Value = "中庚价值灵动灵活配置混合"
$newQ4.Range("D4").Value = "34.74"
$newQ4.Range("E4").Value = "93.96"
$newQ4.Range("F4").Value = "4.83"
$newQ4.Range("G4").Value = "1.6779"
$newQ4.Range("H4").Value = 2
$newQ4.Range("A5").Value = 3
$newQ4.Range("B5").Value = "481010"
$newQ4.Range("C5").Value = "工银中小盘混合"
$newQ4.Range("D5").Value = "15.71"
$newQ4.Range("E5").Value = "91.30"
$newQ4.Range("F5").Value = "4.87"
$newQ4.Range("G5").Value = "0.7651"
$newQ4.Range("H5").Value = 1
$newQ4.Range("A6").Value = 4
$newQ4.Range("B6").Value = "005939"
$newQ4.Range("C6").Value = "工银新能源汽车混合A"
$newQ4.Range("D6").Value = "30.95"
$newQ4.Range("E6").Value = "81.75"
$newQ4.Range("F6").Value = "2.00"
$newQ4.Range("G6").Value = "0.6190"
$newQ4.Range("H6").Value = 9
$newQ4.Range("A7").Value = 5
$newQ4.Range("B7").Value = "005940"
$newQ4.Range("C7").Value = "工银新能源汽车混合C"
$newQ4.Range("D7").Value = "26.09"
$newQ4.Range("E7").Value = "81.75"
$newQ4.Range("F7").Value = "2.00"
$newQ4.Range("G7").Value = "0.5218"
$newQ4.Range("H7").Value = 9
$newQ4.Range("A8").Value = 6
$newQ4.Range("B8").Value = "100029"
$newQ4.Range("C8").Value = "富国天成红利混合"
$newQ4.Range("D8").Value = "8.86"
$newQ4.Range("E8").Value = "73.21"
$newQ4.Range("F8").Value = "2.68"
$newQ4.Range("G8").Value = "0.2374"
$newQ4.Range("H8").Value = 6
$newQ4.Range("A9").Value = 7
$newQ4.Range("B9").Value = "013533"
$newQ4.Range("C9").Value = "广发科技创新混合C"
$newQ4.Range("D9").Value = "2.90"
$newQ4.Range("E9").Value = "91.92"
$newQ4.Range("F9").Value = "6.45"
$newQ4.Range("G9").Value = "0.1870"
$newQ4.Range("H9").Value = 3
$newQ4.Range("A10").Value = 8
$newQ4.Range("B10").Value = "373010"
$newQ4.Range("C10").Value = "上投摩根双息平衡混合A"
$newQ4.Range("D10").Value = "8.17"
$newQ4.Range("E10").Value = "59.37"
$newQ4.Range("F10").Value = "2.17"
$newQ4.Range("G10").Value = "0.1773"
$newQ4.Range("H10").Value = 6
$newQ4.Range("A11").Value = 9
$newQ4.Range("B11").Value = "001716"
$newQ4.Range("C11").Value = "工银新趋势灵活配置混合A"
$newQ4.Range("D11").Value = "4.17"
$newQ4.Range("E11").Value = "81.40"
$newQ4.Range("F11").Value = "4.05"
$newQ4.Range("G11").Value = "0.1689"
$newQ4.Range("H11").Value = 5
$newQ4.Range("A12").Value = 10
$newQ4.Range("B12").Value = "002810"
$newQ4.Range("C12").Value = "金信转型创新成长灵活配置混合"
$newQ4.Range("D12").Value = "3.84"
$newQ4.Range("E12").Value = "89.18"
$newQ4.Range("F12").Value = "3.09"
$newQ4.Range("G12").Value = "0.1187"
$newQ4.Range("H12").Value = 9
$newQ4.Range("A13").Value = 11
$newQ4.Range("B13").Value = "009874"
$newQ4.Range("C13").Value = "九泰久睿量化股票A"
$newQ4.Range("D13").Value = "2.60"
$newQ4.Range("E13").Value = "92.28"
$newQ4.Range("F13").Value = "3.06"
$newQ4.Range("G13").Value = "0.0796"
$newQ4.Range("H13").Value = 7
$newQ4.Range("A14").Value = 12
$newQ4.Range("B14").Value = "010779"
$newQ4.Range("C14").Value = "西部利得量化优选一年持有期混合A"
$newQ4.Range("D14").Value = "4.22"
$newQ4.Range("E14").Value = "89.52"
$newQ4.Range("F14").Value = "1.82"
$newQ4.Range("G14").Value = "0.0768"
$newQ4.Range("H14").Value = 9
$newQ4.Range("A15").Value = 13
$newQ4.Range("B15").Value = "001997"
$newQ4.Range("C15").Value = "工银新趋势灵活配置混合C"
$newQ4.Range("D15").Value = "1.53"
$newQ4.Range("E15").Value = "81.40"
$newQ4.Range("F15").Value = "4.05"
$newQ4.Range("G15").Value = "0.0620"
$newQ4.Range("H15").Value = 5
$newQ4.Range("A16").Value = 14
$newQ4.Range("B16").Value = "590006"
$newQ4.Range("C16").Value = "中邮中小盘灵活配置混合"
$newQ4.Range("D16").Value = "2.56"
$newQ4.Range("E16").Value = "74.40"
$newQ4.Range("F16").Value = "2.08"
$newQ4.Range("G16").Value = "0.0532"
$newQ4.Range("H16").Value = 6
$newQ4.Range("A17").Value = 15
$newQ4.Range("B17").Value = "010780"
$newQ4.Range("C17").Value = "西部利得量化优选一年持有期混合C"
$newQ4.Range("D17").Value = "1.88"
$newQ4.Range("E17").Value = "89.52"
$newQ4.Range("F17").Value = "1.82"
$newQ4.Range("G17").Value = "0.0342"
$newQ4.Range("H17").Value = 9
$newQ4.Range("A18").Value = 16
$newQ4.Range("B18").Value = "002485"
$newQ4.Range("C18").Value = "国联安通盈灵活配置混合C"
$newQ4.Range("D18").Value = "1.75"
$newQ4.Range("E18").Value = "36.18"
$newQ4.Range("F18").Value = "1.23"
$newQ4.Range("G18").Value = "0.0215"
$newQ4.Range("H18").Value = 8
$newQ4.Range("A19").Value = 17
$newQ4.Range("B19").Value = "519139"
$newQ4.Range("C19").Value = "海富通沪港深灵活配置混合"
$newQ4.Range("D19").Value = "0.67"
$newQ4.Range("E19").Value = "92.35"
$newQ4.Range("F19").Value = "3.17"
$newQ4.Range("G19").Value = "0.0212"
$newQ4.Range("H19").Value = 6
$newQ4.Range("A20").Value = 18
$newQ4.Range("B20").Value = "011500"
$newQ4.Range("C20").Value = "九泰量化新兴产业混合"
$newQ4.Range("D20").Value = "0.50"
$newQ4.Range("E20").Value = "91.24"
$newQ4.Range("F20").Value = "3.74"
$newQ4.Range("G20").Value = "0.0187"
$newQ4.Range("H20").Value = 3
$newQ4.Range("A21").Value = 19
$newQ4.Range("B21").Value = "001897"
$newQ4.Range("C21").Value = "九泰久盛量化先锋灵活配置混合A"
$newQ4.Range("D21").Value = "0.45"
$newQ4.Range("E21").Value = "92.70"
$newQ4.Range("F21").Value = "3.52"
$newQ4.Range("G21").Value = "0.0158"
$newQ4.Range("H21").Value = 8
$newQ4.Range("A22").Value = 20
$newQ4.Range("B22").Value = "010703"
$newQ4.Range("C22").Value = "财通智选消费股票A"
$newQ4.Range("D22").Value = "0.45"
$newQ4.Range("E22").Value = "92.71"
$newQ4.Range("F22").Value = "3.10"
$newQ4.Range("G22").Value = "0.0140"
$newQ4.Range("H22").Value = 7
$newQ4.Range("A23").Value = 21
$newQ4.Range("B23").Value = "002186"
$newQ4.Range("C23").Value = "国联安鑫享灵活配置混合C"
$newQ4.Range("D23").Value = "0.84"
$newQ4.Range("E23").Value = "30.88"
$newQ4.Range("F23").Value = "1.63"
$newQ4.Range("G23").Value = "0.0137"
$newQ4.Range("H23").Value = 5
$newQ4.Range("A24").Value = 22
$newQ4.Range("B24").Value = "010120"
$newQ4.Range("C24").Value = "九泰久福量化股票A"
$newQ4.Range("D24").Value = "0.44"
$newQ4.Range("E24").Value = "93.44"
$newQ4.Range("F24").Value = "3.10"
$newQ4.Range("G24").Value = "0.0136"
$newQ4.Range("H24").Value = 8
$newQ4.Range("A25").Value = 23
$newQ4.Range("B25").Value = "010704"
$newQ4.Range("C25").Value = "财通智选消费股票C"
$newQ4.Range("D25").Value = "0.44"
$newQ4.Range("E25").Value = "92.71"
$newQ4.Range("F25").Value = "3.10"
$newQ4.Range("G25").Value = "0.0136"
$newQ4.Range("H25").Value = 7
$newQ4.Range("A26").Value = 24
$newQ4.Range("B26").Value = "009043"
$newQ4.Range("C26").Value = "九泰久信量化股票"
$newQ4.Range("D26").Value = "0.37"
$newQ4.Range("E26").Value = "92.67"
$newQ4.Range("F26").Value = "3.02"
$newQ4.Range("G26").Value = "0.0112"
$newQ4.Range("H26").Value = 8
$newQ4.Range("A27").Value = 25
$newQ4.Range("B27").Value = "007903"
$newQ4.Range("C27").Value = "长城量化小盘股票"
$newQ4.Range("D27").Value = "0.84"
$newQ4.Range("E27").Value = "89.96"
$newQ4.Range("F27").Value = "1.14"
$newQ4.Range("G27").Value = "0.0096"
$newQ4.Range("H27").Value = 5
$newQ4.Range("A28").Value = 26
$newQ4.Range("B28").Value = "000664"
$newQ4.Range("C28").Value = "国联安通盈灵活配置混合A"
$newQ4.Range("D28").Value = "0.68"
$newQ4.Range("E28").Value = "36.18"
$newQ4.Range("F28").Value = "1.23"
$newQ4.Range("G28").Value = "0.0084"
$newQ4.Range("H28").Value = 8
$newQ4.Range("A29").Value = 27
$newQ4.Range("B29").Value = "013242"
$newQ4.Range("C29").Value = "北信瑞丰优势行业股票"
$newQ4.Range("D29").Value = "0.49"
$newQ4.Range("E29").Value = "91.56"
$newQ4.Range("F29").Value = "1.71"
$newQ4.Range("G29").Value = "0.0084"
$newQ4.Range("H29").Value = 8
$newQ4.Range("A30").Value = 28
$newQ4.Range("B30").Value = "007527"
$newQ4.Range("C30").Value = "融通量化多策略灵活配置混合A"
$newQ4.Range("D30").Value = "0.28"
$newQ4.Range("E30").Value = "91.43"
$newQ4.Range("F30").Value = "2.51"
$newQ4.Range("G30").Value = "0.0070"
$newQ4.Range("H30").Value = 4
$newQ4.Range("A31").Value = 29
$newQ4.Range("B31").Value = "015633"
$newQ4.Range("C31").Value = "中金景气驱动混合A"
$newQ4.Range("D31").Value = "0.09"
$newQ4.Range("E31").Value = "90.15"
$newQ4.Range("F31").Value = "3.18"
$newQ4.Range("G31").Value = "0.0029"
$newQ4.Range("H31").Value = 6
$newQ4.Range("A32").Value = 30
$newQ4.Range("B32").Value = "009054"
$newQ4.Range("C32").Value = "圆信永丰沣泰混合"
$newQ4.Range("D32").Value = "0.23"
$newQ4.Range("E32").Value = "31.90"
$newQ4.Range("F32").Value = "1.13"
$newQ4.Range("G32").Value = "0.0026"
$newQ4.Range("H32").Value = 9
$newQ4.Range("A33").Value = 31
$newQ4.Range("B33").Value = "015634"
$newQ4.Range("C33").Value = "中金景气驱动混合C"
$newQ4.Range("D33").Value = "0.06"
$newQ4.Range("E33").Value = "90.15"
$newQ4.Range("F33").Value = "3.18"
$newQ4.Range("G33").Value = "0.0019"
$newQ4.Range("H33").Value = 6
$newQ4.Range("A34").Value = 32
$newQ4.Range("B34").Value = "004510"
$newQ4.Range("C34").Value = "九泰久盛量化先锋灵活配置混合C"
$newQ4.Range("D34").Value = "0.04"
$newQ4.Range("E34").Value = "92.70"
$newQ4.Range("F34").Value = "3.52"
$newQ4.Range("G34").Value = "0.0014"
$newQ4.Range("H34").Value = 8
$newQ4.Range("A35").Value = 33
$newQ4.Range("B35").Value = "001228"
$newQ4.Range("C35").Value = "国联安鑫享灵活配置混合A"
$newQ4.Range("D35").Value = "0.08"
$newQ4.Range("E35").Value = "30.88"
$newQ4.Range("F35").Value = "1.63"
$newQ4.Range("G35").Value = "0.0013"
$newQ4.Range("H35").Value = 5
$newQ4.Range("A36").Value = 34
$newQ4.Range("B36").Value = "010121"
$newQ4.Range("C36").Value = "九泰久福量化股票C"
$newQ4.Range("D36").Value = "0.03"
$newQ4.Range("E36").Value = "93.44"
$newQ4.Range("F36").Value = "3.10"
$newQ4.Range("G36").Value = "0.0009"
$newQ4.Range("H36").Value = 8
$newQ4.Range("A37").Value = 35
$newQ4.Range("B37").Value = "007528"
$newQ4.Range("C37").Value = "融通量化多策略灵活配置混合C"
$newQ4.Range("D37").Value = "0.03"
$newQ4.Range("E37").Value = "91.43"
$newQ4.Range("F37").Value = "2.51"
$newQ4.Range("G37").Value = "0.0008"
$newQ4.Range("H37").Value = 4
$newQ4.Range("A38").Value = 36
$newQ4.Range("B38").Value = "960005"
$newQ4.Range("C38").Value = "上投摩根双息平衡混合H"
$newQ4.Range("D38").Value = "0.02"
$newQ4.Range("E38").Value = "59.37"
$newQ4.Range("F38").Value = "2.17"
$newQ4.Range("G38").Value = "0.0004"
$newQ4.Range("H38").Value = 6
$newQ4.Range("A39").Value = 37
$newQ4.Range("B39").Value = "016803"
$newQ4.Range("C39").Value = "上投摩根双息平衡混合C"
$newQ4.Range("D39").Value = "0.00"
$newQ4.Range("E39").Value = "59.37"
$newQ4.Range("F39").Value = "2.17"
$newQ4.Range("G39").Value = 0
$newQ4.Range("H39").Value = 6
$newQ4.Range("A40").Value = 38
$newQ4.Range("B40").Value = "016399"
$newQ4.Range("C40").Value = "九泰久睿量化股票C"
$newQ4.Range("D40").Value = "0.00"
$newQ4.Range("E40").Value = "92.28"
$newQ4.Range("F40").Value = "3.06"
$newQ4.Range("G40").Value = 0
$newQ4.Range("H40").Value = 7

# --- Step 4: update the "总计" (summary) sheet: insert 2022-Q4 as new row 2, shift the rest down ---
$zj = $wb.Worksheets.Item("总计")

# create row 9 (new) by copying the style of row 8 (column A) first
$zj.Range("A8").Copy($zj.Range("A9"))

$zj.Range("A2").Value = 0
$zj.Range("B2").Value = "2022-Q4"
$zj.Range("C2").Value = 39
$zj.Range("D2").Value = 10.97
$zj.Range("A3").Value = 1
$zj.Range("B3").Value = "2022-Q3"
$zj.Range("C3").Value = 33
$zj.Range("D3").Value = 9.41
$zj.Range("A4").Value = 2
$zj.Range("B4").Value = "2022-Q2"
$zj.Range("C4").Value = 11
$zj.Range("D4").Value = 3.74
$zj.Range("A5").Value = 3
$zj.Range("B5").Value = "2022-Q1"
$zj.Range("C5").Value = 2
$zj.Range("D5").Value = 0.23
$zj.Range("A6").Value = 4
$zj.Range("B6").Value = "2021-Q4"
$zj.Range("C6").Value = 7
$zj.Range("D6").Value = 0.33
$zj.Range("A7").Value = 5
$zj.Range("B7").Value = "2021-Q2"
$zj.Range("C7").Value = 8
$zj.Range("D7").Value = 0.2
$zj.Range("A8").Value = 6
$zj.Range("B8").Value = "2021-Q1"
$zj.Range("C8").Value = 4
$zj.Range("D8").Value = 0.07
$zj.Range("A9").Value = 7
$zj.Range("B9").Value = "2020-Q4"
$zj.Range("C9").Value = 1
$zj.Range("D9").Value = 0.01

Write-Host "done"
